$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old data rows (rows 2-6), then rewrite the table with the
# new item list (Xoài lắc / Nước mía) replacing the old one.
$ws.Range("A2:D6").ClearContents()

$ws.Range("A2").Value = "Xoài lắc"
$ws.Range("B2").Value = "Dĩa"
$ws.Range("C2").Value = 20000
$ws.Range("D2").Value = 22000

$ws.Range("A3").Value = "Nước mía"
$ws.Range("B3").Value = "Ly"
$ws.Range("C3").Value = 3000
$ws.Range("D3").Value = 6000

# Remove now-unused rows 4-6 entirely so dimension shrinks to A1:D3
$ws.Range("A4:D6").Delete()

$ws.Range("A4").Select()
